$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "65.994.38"
$ws.Range("E2").Value = "  +1.65%  "

# Row 3
$ws.Range("D3").Value = "3.210.19"
$ws.Range("E3").Value = "  +1.24%  "

# Row 4
$ws.Range("E4").Value = "  -0.13%  "

# Row 5
$ws.Range("D5").Value = "606.11"
$ws.Range("E5").Value = "  +4.48%  "

# Row 6
$ws.Range("D6").Value = "153.04"
$ws.Range("E6").Value = "  +0.74%  "

# Row 7
$ws.Range("D7").Value = "0.999"
$ws.Range("E7").Value = "  -0.05%  "

# Row 8
$ws.Range("D8").Value = "3.208.47"
$ws.Range("E8").Value = "  +1.23%  "

# Row 9
$ws.Range("E9").Value = "  -0.09%  "

# Row 10
$ws.Range("E10").Value = "  -0.99%  "

# Row 11
$ws.Range("D11").Value = "6.17"
$ws.Range("E11").Value = "  -1.21%  "

# Row 12
$ws.Range("D12").Value = "0.507"
$ws.Range("E12").Value = "  +1.01%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000271"
$ws.Range("E13").Value = "  +0.42%  "

# Row 14
$ws.Range("D14").Value = "38.98"
$ws.Range("E14").Value = "  +3.53%  "

# Row 15
$ws.Range("D15").Value = "3.735.93"
$ws.Range("E15").Value = "  +1.12%  "

# Row 16
$ws.Range("D16").Value = "66.106.76"
$ws.Range("E16").Value = "  +1.63%  "

# Row 17
$ws.Range("D17").Value = "7.46"
$ws.Range("E17").Value = "  +3.84%  "

# Row 18
$ws.Range("D18").Value = "3.224.92"
$ws.Range("E18").Value = "  +1.75%  "

# Row 19
$ws.Range("E19").Value = "  -0.32%  "

# Row 20
$ws.Range("D20").Value = "510.16"
$ws.Range("E20").Value = "  -0.66%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "15.50"
$ws.Range("E21").Value = "  +4.50%  "

# Row 22
$ws.Range("D22").Value = "0.733"
$ws.Range("E22").Value = "  +0.59%  "

# Row 23
$ws.Range("D23").Value = "15.35"
$ws.Range("E23").Value = "  +0.63%  "

# Row 24
$ws.Range("D24").Value = "8.03"
$ws.Range("E24").Value = "  +2.98%  "

# Row 25
$ws.Range("D25").Value = "85.17"
$ws.Range("E25").Value = "  -0.22%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.00"
$ws.Range("E26").Value = "  +0.11%  "

# Row 27
$ws.Range("D27").Value = "3.02"
$ws.Range("E27").Value = "  +2.91%  "

# Row 28
$ws.Range("D28").Value = "9.12"
$ws.Range("E28").Value = "  +1.15%  "

# Row 29
$ws.Range("D29").Value = "2.24"
$ws.Range("E29").Value = "  +2.62%  "

# Row 30
$ws.Range("E30").Value = "  +2.61%  "

# Row 31
$ws.Range("B31").Value = "NEARProtocol"
$ws.Range("C31").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D31").Value = "6.78"
$ws.Range("E31").Value = "  +6.93%  "

# Row 32
$ws.Range("B32").Value = "EthereumClassic"
$ws.Range("C32").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D32").Value = "28.08"
$ws.Range("E32").Value = "  +0.72%  "

# Row 33
$ws.Range("B33").Value = "Mantle"
$ws.Range("C33").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D33").Value = "1.21"
$ws.Range("E33").Value = "  +1.16%  "

# Row 34
$ws.Range("B34").Value = "FirstDigitalUSD"
$ws.Range("C34").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.00"
$ws.Range("E34").Value = "  +0.05%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "6.60"
$ws.Range("E35").Value = "  +0.41%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "55.30"
$ws.Range("E36").Value = "  -0.82%  "

# Row 37
$ws.Range("D37").Value = "0.0903"
$ws.Range("E37").Value = "  +1.01%  "

# Row 38
$ws.Range("D38").Value = "480.87"
$ws.Range("E38").Value = "  +1.20%  "

# Row 39
$ws.Range("D39").Value = "0.0419"
$ws.Range("E39").Value = "  -0.19%  "

# Row 40
$ws.Range("D40").Value = "2.96"
$ws.Range("E40").Value = "  -6.34%  "

# Row 41
$ws.Range("D41").Value = "8.85"
$ws.Range("E41").Value = "  +2.24%  "

# Row 42
$ws.Range("D42").Value = "0.296"
$ws.Range("E42").Value = "  +3.70%  "

# Row 43
$ws.Range("E43").Value = "  -0.29%  "

# Row 44
$ws.Range("D44").Value = "2.936.05"
$ws.Range("E44").Value = "  -4.16%  "

# Row 45
$ws.Range("D45").Value = "2.45"
$ws.Range("E45").Value = "  +1.05%  "

# Row 46
$ws.Range("D46").Value = "0.0₃0640"
$ws.Range("E46").Value = "  +4.51%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "28.60"

# Row 48
$ws.Range("E48").Value = "  +0.06%  "

# Row 49
$ws.Range("D49").Value = "0.116"
$ws.Range("E49").Value = "  +0.27%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.30"
$ws.Range("E50").Value = "  +2.05%  "

# Row 51
$ws.Range("B51").Value = "Monero"
$ws.Range("C51").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D51").Value = "120.29"
$ws.Range("E51").Value = "  -0.21%  "
